$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8958579383065057
$ws.Range("C2").Value = 0.2247595374067828
$ws.Range("D2").Value = 0.3319435897582679
$ws.Range("F2").Value = 1.174130829997679
$ws.Range("G2").Value = 0.002417822349252524
$ws.Range("I2").Value = 0.4468222055563444
$ws.Range("J2").Value = 0.3476189276469199
$ws.Range("N2").Value = 0.9778347057375498
$ws.Range("O2").Value = 2.412498572457594

$ws.Range("B3").Value = 0.7985604697435065
$ws.Range("C3").Value = 0.1971300892429326
$ws.Range("D3").Value = 0.3232883474917969
$ws.Range("F3").Value = 1.166228234281647
$ws.Range("G3").Value = 0.00242073139868274
$ws.Range("I3").Value = 0.4518980869209024
$ws.Range("J3").Value = 0.3360577617490037
$ws.Range("N3").Value = 0.9814234588060984
$ws.Range("O3").Value = 2.407310167648546

$ws.Range("B4").Value = 0.7388045868442532
$ws.Range("C4").Value = 0.1801297771864654
$ws.Range("D4").Value = 0.3181145741445448
$ws.Range("F4").Value = 1.162158867636791
$ws.Range("G4").Value = 0.002422613250547532
$ws.Range("I4").Value = 0.4553437244751066
$ws.Range("J4").Value = 0.329167463114203
$ws.Range("N4").Value = 0.9840343465379391
$ws.Range("O4").Value = 2.405881948729188

$ws.Range("B5").Value = 0.7144512798417111
$ws.Range("C5").Value = 0.1731933604411608
$ws.Range("D5").Value = 0.3160416689919856
$ws.Range("F5").Value = 1.160697185313637
$ws.Range("G5").Value = 0.00242340425413124
$ws.Range("I5").Value = 0.4568304185298011
$ws.Range("J5").Value = 0.3264119005761472
$ws.Range("N5").Value = 0.9852009196453864
$ws.Range("O5").Value = 2.405741073243689

$ws.Range("B6").Value = 0.710407335169549
$ws.Range("C6").Value = 0.1720410609483451
$ws.Range("D6").Value = 0.3156996086723893
$ws.Range("F6").Value = 1.16046634235343
$ws.Range("G6").Value = 0.002423537059486641
$ws.Range("I6").Value = 0.4570822649938364
$ws.Range("J6").Value = 0.3259574971127677
$ws.Range("N6").Value = 0.9854008312185911
$ws.Range("O6").Value = 2.405744301159444

$ws.Range("B7").Value = 0.7384761570394005
$ws.Range("C7").Value = 0.1800362648334328
$ws.Range("D7").Value = 0.3180864745549883
$ws.Range("F7").Value = 1.162138359099188
$ws.Range("G7").Value = 0.002422623820549365
$ws.Range("I7").Value = 0.455363440456523
$ws.Range("J7").Value = 0.3291300890307269
$ws.Range("N7").Value = 0.9840496636299534
$ws.Range("O7").Value = 2.405878263757614

$ws.Range("B8").Value = 0.8623136804514502
$ws.Range("C8").Value = 0.2152405067602103
$ws.Range("D8").Value = 0.3289301409490264
$ws.Range("F8").Value = 1.171243382787537
$ws.Range("G8").Value = 0.002418805574423825
$ws.Range("I8").Value = 0.4485039794284411
$ws.Range("J8").Value = 0.3435893457580477
$ws.Range("N8").Value = 0.9789876592773723
$ws.Range("O8").Value = 2.410344355260435

$ws.Range("B9").Value = 1.104989903737021
$ws.Range("C9").Value = 0.2839813169477168
$ws.Range("D9").Value = 0.3513070242476886
$ws.Range("F9").Value = 1.195323465710231
$ws.Range("G9").Value = 0.002412073867122722
$ws.Range("I9").Value = 0.43767087619036
$ws.Range("J9").Value = 0.373602642138934
$ws.Range("N9").Value = 0.9722865148499977
$ws.Range("O9").Value = 2.433087569356871

$ws.Range("B10").Value = 1.283130306744283
$ws.Range("C10").Value = 0.3342954703116447
$ws.Range("D10").Value = 0.3684236380674122
$ws.Range("F10").Value = 1.216832826748004
$ws.Range("G10").Value = 0.002407584185763126
$ws.Range("I10").Value = 0.4313185726825637
$ws.Range("J10").Value = 0.3966755800465904
$ws.Range("N10").Value = 0.9693213624961885
$ws.Range("O10").Value = 2.458385010120224

$ws.Range("B11").Value = 1.36412791747756
$ws.Range("C11").Value = 0.3571416509150254
$ws.Range("D11").Value = 0.3763570209600857
$ws.Range("F11").Value = 1.227452118134494
$ws.Range("G11").Value = 0.002405639754356148
$ws.Range("I11").Value = 0.4287799373963388
$ws.Range("J11").Value = 0.4073965349386697
$ws.Range("N11").Value = 0.9683960948349721
$ws.Range("O11").Value = 2.47177208133823

$ws.Range("B12").Value = 1.39479274039553
$ws.Range("C12").Value = 0.3657866028133299
$ws.Range("D12").Value = 0.3793822498334976
$ws.Range("F12").Value = 1.231593716365438
$ws.Range("G12").Value = 0.002404917457782128
$ws.Range("I12").Value = 0.427869293296137
$ws.Range("J12").Value = 0.4114887751559735
$ws.Range("N12").Value = 0.968106501462529
$ws.Range("O12").Value = 2.477112633217587

$ws.Range("B13").Value = 1.388188867690701
$ws.Range("C13").Value = 0.3639250482097509
$ws.Range("D13").Value = 0.3787297790239847
$ws.Range("F13").Value = 1.230696393918379
$ws.Range("G13").Value = 0.002405072395041125
$ws.Range("I13").Value = 0.4280631595147781
$ws.Range("J13").Value = 0.4106059933379385
$ws.Range("N13").Value = 0.9681661690905088
$ws.Range("O13").Value = 2.475950374984535

$ws.Range("B14").Value = 1.366650884393721
$ws.Range("C14").Value = 0.357853006674361
$ws.Range("D14").Value = 0.3766054874928955
$ws.Range("F14").Value = 1.227790437187267
$ws.Range("G14").Value = 0.002405580050034653
$ws.Range("I14").Value = 0.4287040012694519
$ws.Range("J14").Value = 0.407732555148911
$ws.Range("N14").Value = 0.9683710522911468
$ws.Range("O14").Value = 2.472206011946128

$ws.Range("B15").Value = 1.353457262123072
$ws.Range("C15").Value = 0.354132861539199
$ws.Range("D15").Value = 0.3753070331537742
$ws.Range("F15").Value = 1.226026130463623
$ws.Range("G15").Value = 0.002405892826963068
$ws.Range("I15").Value = 0.4291031415777624
$ws.Range("J15").Value = 0.4059767198024247
$ws.Range("N15").Value = 0.9685044615502107
$ws.Range("O15").Value = 2.469947822234587

$ws.Range("B16").Value = 1.277835997430202
$ws.Range("C16").Value = 0.3328015427737512
$ws.Range("D16").Value = 0.3679081210939898
$ws.Range("F16").Value = 1.216155649838655
$ws.Range("G16").Value = 0.002407713224163016
$ws.Range("I16").Value = 0.431491559782792
$ws.Range("J16").Value = 0.3959794756007966
$ws.Range("N16").Value = 0.9693903437363929
$ws.Range("O16").Value = 2.457548027018561

$ws.Range("B17").Value = 1.231433613381398
$ws.Range("C17").Value = 0.3197044579949022
$ws.Range("D17").Value = 0.3634066926782964
$ws.Range("F17").Value = 1.210314392236171
$ws.Range("G17").Value = 0.002408855017714254
$ws.Range("I17").Value = 0.4330468340654932
$ws.Range("J17").Value = 0.3899041776225545
$ws.Range("N17").Value = 0.9700422043823949
$ws.Range("O17").Value = 2.450423116093674

$ws.Range("B18").Value = 1.204740578015958
$ws.Range("C18").Value = 0.3121674315626137
$ws.Range("D18").Value = 0.3608314277825855
$ws.Range("F18").Value = 1.207033192039489
$ws.Range("G18").Value = 0.002409520970268854
$ws.Range("I18").Value = 0.4339744181033751
$ws.Range("J18").Value = 0.3864309953096949
$ws.Range("N18").Value = 0.9704570153686944
$ws.Range("O18").Value = 2.446501872926746

$ws.Range("B19").Value = 1.19570220212546
$ws.Range("C19").Value = 0.3096148601057678
$ws.Range("D19").Value = 0.3599618670883444
$ws.Range("F19").Value = 1.205935713887143
$ws.Range("G19").Value = 0.00240974803645635
$ws.Range("I19").Value = 0.4342941492036232
$ws.Range("J19").Value = 0.3852586681516499
$ws.Range("N19").Value = 0.9706043162922384
$ws.Range("O19").Value = 2.445204547230617

$ws.Range("B20").Value = 1.236373611424881
$ws.Range("C20").Value = 0.3210990743784237
$ws.Range("D20").Value = 0.3638844460114967
$ws.Range("F20").Value = 1.210928073709468
$ws.Range("G20").Value = 0.002408732517649161
$ws.Range("I20").Value = 0.432877852385861
$ws.Range("J20").Value = 0.3905487119038327
$ws.Range("N20").Value = 0.9699686864184685
$ws.Range("O20").Value = 2.451163268291282

$ws.Range("B21").Value = 1.372977320344035
$ws.Range("C21").Value = 0.3596366890755007
$ws.Range("D21").Value = 0.377228873503185
$ws.Range("F21").Value = 1.228640719801049
$ws.Range("G21").Value = 0.002405430559557887
$ws.Range("I21").Value = 0.428514393422418
$ws.Range("J21").Value = 0.4085756722045915
$ws.Range("N21").Value = 0.9683092244138294
$ws.Range("O21").Value = 2.4732984554729

$ws.Range("B22").Value = 1.462212971771066
$ws.Range("C22").Value = 0.3847857708007041
$ws.Range("D22").Value = 0.3860727538397839
$ws.Range("F22").Value = 1.240918311918833
$ws.Range("G22").Value = 0.002403354213214722
$ws.Range("I22").Value = 0.4259581270551109
$ws.Range("J22").Value = 0.4205465061179865
$ws.Range("N22").Value = 0.9675789190679609
$ws.Range("O22").Value = 2.489346036005969

$ws.Range("B23").Value = 1.414590651192782
$ws.Range("C23").Value = 0.3713667817498276
$ws.Range("D23").Value = 0.3813414310696999
$ws.Range("F23").Value = 1.234301259432939
$ws.Range("G23").Value = 0.002404454947170774
$ws.Range("I23").Value = 0.4272953495551199
$ws.Range("J23").Value = 0.4141401048818238
$ws.Range("N23").Value = 0.9679363234536851
$ws.Range("O23").Value = 2.480636171849511

$ws.Range("B24").Value = 1.234140287051957
$ws.Range("C24").Value = 0.3204685911438503
$ws.Range("D24").Value = 0.3636684142691138
$ws.Range("F24").Value = 1.210650388432398
$ws.Range("G24").Value = 0.002408787870184122
$ws.Range("I24").Value = 0.432954144864734
$ws.Range("J24").Value = 0.3902572569538307
$ws.Range("N24").Value = 0.9700017991312961
$ws.Range("O24").Value = 2.450828100544925

$ws.Range("B25").Value = 1.039362909727629
$ws.Range("C25").Value = 0.2654176967645583
$ws.Range("D25").Value = 0.3451345807787618
$ws.Range("F25").Value = 1.188140374210093
$ws.Range("G25").Value = 0.002413814538741788
$ws.Range("I25").Value = 0.4403201174725062
$ws.Range("J25").Value = 0.3653045385717206
$ws.Range("N25").Value = 0.9737549109638337
$ws.Range("O25").Value = 2.42543128245191
